$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()

$ws.Range("A8").Value = 106
$ws.Range("C8").Value = "暂时移除角色"
$ws.Range("D8").Value = "[[107:1003]]"

$ws.Range("A9").Value = 107
$ws.Range("C9").Value = "加回角色"
$ws.Range("D9").Value = "[[108:1003:26]]"

$ws.Range("B8").Value = "Option6"
$ws.Range("B9").Value = "Option7"

$ws.Rows.Item(8).Select() | Out-Null
